# Implement change requests for forms
# - survey sheet: wrap hh_death_id/name/surname block in begin/end screen,
#   add a second begin screen before the gender/note/date block, and shift
#   rows accordingly (A1:L16 -> A1:L18)
# - settings sheet: bump form_version
# - active sheet/selection moves from survey to settings

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rebuild the "survey" sheet content with the new row layout.
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Cells.Clear()

# Row 1 - header (unchanged)
$survey.Cells.Item(1, 1).Value = 'clause'
$survey.Cells.Item(1, 2).Value = 'condition'
$survey.Cells.Item(1, 3).Value = 'type'
$survey.Cells.Item(1, 4).Value = 'values_list'
$survey.Cells.Item(1, 5).Value = 'name'
$survey.Cells.Item(1, 6).Value = 'display.prompt'
$survey.Cells.Item(1, 7).Value = 'inputAttributes.min'
$survey.Cells.Item(1, 8).Value = 'inputAttributes.step'
$survey.Cells.Item(1, 9).Value = 'calculation'
$survey.Cells.Item(1, 10).Value = 'constraint'
$survey.Cells.Item(1, 11).Value = 'display.constraint_message'
$survey.Cells.Item(1, 12).Value = 'hideInContents'

# Row 2 - hh_death_id question (moved up, used to be row 3)
$survey.Cells.Item(2, 3).Value = 'text'
$survey.Cells.Item(2, 5).Value = 'hh_death_id'
$survey.Cells.Item(2, 6).Value = 'q65a'
$survey.Cells.Item(2, 10).Value = "/^[A-Z]{3}-[0-9]{3}-7[0-9]{2}$/.test(data('hh_death_id'))"
$survey.Cells.Item(2, 11).Value = 'extid_format'

# Row 3 - NEW: begin screen (wraps id/name/surname block)
$survey.Cells.Item(3, 1).Value = 'begin screen'

# Row 4 - hh_death_name question (unchanged position)
$survey.Cells.Item(4, 3).Value = 'text'
$survey.Cells.Item(4, 5).Value = 'hh_death_name'
$survey.Cells.Item(4, 6).Value = 'q65b'
$survey.Cells.Item(4, 10).Value = "!/\p{N}/u.test(data('hh_death_name'))"
$survey.Cells.Item(4, 11).Value = 'name_number'

# Row 5 - hh_death_surname question (unchanged position)
$survey.Cells.Item(5, 3).Value = 'text'
$survey.Cells.Item(5, 5).Value = 'hh_death_surname'
$survey.Cells.Item(5, 6).Value = 'q65c'
$survey.Cells.Item(5, 10).Value = "!/\p{N}/u.test(data('hh_death_surname'))"
$survey.Cells.Item(5, 11).Value = 'surname_number'

# Row 6 - NEW: end screen (closes id/name/surname block)
$survey.Cells.Item(6, 1).Value = 'end screen'

# Row 7 - gender question (used to be row 6)
$survey.Cells.Item(7, 3).Value = 'select_one'
$survey.Cells.Item(7, 4).Value = 'gender'
$survey.Cells.Item(7, 5).Value = 'hh_death_gender'
$survey.Cells.Item(7, 6).Value = 'q65d'

# Row 8 - NEW: begin screen (opens gender/date block)
$survey.Cells.Item(8, 1).Value = 'begin screen'

# Row 9 - note question (used to be row 7)
$survey.Cells.Item(9, 3).Value = 'note'
$survey.Cells.Item(9, 6).Value = 'q65e'

# Row 10 - if (used to be row 8)
$survey.Cells.Item(10, 1).Value = 'if'
$survey.Cells.Item(10, 2).Value = "not(selected(data('hh_death_date_dk'), 'dk'))"

# Row 11 - birth_date (used to be row 9)
$survey.Cells.Item(11, 3).Value = 'birth_date'
$survey.Cells.Item(11, 5).Value = 'hh_death_date'
$survey.Cells.Item(11, 10).Value = "selected(data('hh_death_date_dk'), 'dk') || !data('hh_death_date') || (!data('hh_death_date').isBefore('2020-03-01') && !data('hh_death_date').isAfter())"
$survey.Cells.Item(11, 11).Value = 'invalid_death_date'
$survey.Cells.Item(11, 12).Value = 1

# Row 12 - end if (used to be row 10)
$survey.Cells.Item(12, 1).Value = 'end if'

# Row 13 - select_multiple dk (used to be row 11)
$survey.Cells.Item(13, 3).Value = 'select_multiple'
$survey.Cells.Item(13, 4).Value = 'dk'
$survey.Cells.Item(13, 5).Value = 'hh_death_date_dk'
$survey.Cells.Item(13, 12).Value = 1

# Row 14 - if (used to be row 12)
$survey.Cells.Item(14, 1).Value = 'if'
$survey.Cells.Item(14, 2).Value = "selected(data('hh_death_date_dk'), 'dk')"

# Row 15 - assign (used to be row 13)
$survey.Cells.Item(15, 3).Value = 'assign'
$survey.Cells.Item(15, 5).Value = 'hh_death_date'
$survey.Cells.Item(15, 9).Value = 'null'

# Row 16 - end if (used to be row 14)
$survey.Cells.Item(16, 1).Value = 'end if'

# Row 17 - NEW: end screen (closes gender/date block)
$survey.Cells.Item(17, 1).Value = 'end screen'

# Row 18 - integer age question (used to be row 15, "end screen" used to be row 16)
$survey.Cells.Item(18, 3).Value = 'integer'
$survey.Cells.Item(18, 5).Value = 'hh_death_age'
$survey.Cells.Item(18, 6).Value = 'q65f'
$survey.Cells.Item(18, 7).Value = 0
$survey.Cells.Item(18, 8).Value = 1

# ---------------------------------------------------------------------
# 2. Bump form_version on the "settings" sheet.
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20210304001

# ---------------------------------------------------------------------
# 3. Move the selection/active sheet: survey -> settings.
# ---------------------------------------------------------------------
$survey.Range("A1").Select()
$settings.Activate()
$settings.Range("B4").Select()
